$d = $word.ActiveDocument
$t = $d.Tables(1)

# Rows 1-3 (1-indexed): "100" -> "0M", "0" -> "0M", "199" -> "0M"
$t.Cell(1, 1).Range.Text = "0M"
$t.Cell(2, 1).Range.Text = "0M"
$t.Cell(3, 1).Range.Text = "0M"

# Insert 10 new rows right after row 3, carrying the new per-iteration values.
# Each insert must re-anchor on the (now shifted) row that used to be row 4, so the
# new rows land in the correct order rather than all stacking before a single anchor.
$newValues = @("103", "0.00003", "0.00007", "0.00004", "0.00000", "0.00004", "0.00004", "0.00004", "0.00421", "100.0")
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $anchor = $t.Rows(4 + $i)
    $newRow = $t.Rows.Add($anchor)
    $t.Cell(4 + $i, 1).Range.Text = $newValues[$i]
}

# The final three rows (previously multi-run, tab-separated summary rows) collapse
# down to a single simple value each.
$n = $t.Rows.Count
$t.Cell($n - 2, 1).Range.Text = "100"
$t.Cell($n - 1, 1).Range.Text = "0"
$t.Cell($n, 1).Range.Text = "199"
